$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: franconian/(196, 11, 189) -> alemannic/(251, 254, 1)
$ws.Range("A7").Value = "(251, 254, 1)"
$ws.Range("B7").Value = "alemannic"

# New rows 13-16
$ws.Range("A13").Value = "(38, 224, 252)"
$ws.Range("B13").Value = "rhine_franconian"

$ws.Range("A14").Value = "(212, 242, 118)"
$ws.Range("B14").Value = "gallo"

$ws.Range("A15").Value = "(244, 191, 162)"
$ws.Range("B15").Value = "gascon"

$ws.Range("A16").Value = "(88, 104, 148)"
$ws.Range("B16").Value = "provencal"

# Column B widens to fit the new, longer culture names (e.g. "rhine_franconian")
$ws.Columns.Item(2).ColumnWidth = 15.45

# Conditional formatting: highlight duplicate values (Home > Conditional Formatting >
# Highlight Cells Rules > Duplicate Values, default red formatting). Rules are added
# A13:A14, B13:B14, A7, B7 (so the sheet XML lists them in that order, matching the
# saved file), then priorities are (re)assigned so B7 is the newest/highest-priority
# rule, matching final numbering 4,3,2,1.
$fcA1314 = $ws.Range("A13:A14").FormatConditions.AddUniqueValues()
$fcA1314.DupeUnique = 1
$fcA1314.Font.Color = 393372
$fcA1314.Interior.Color = 13551615

$fcB1314 = $ws.Range("B13:B14").FormatConditions.AddUniqueValues()
$fcB1314.DupeUnique = 1
$fcB1314.Font.Color = 393372
$fcB1314.Interior.Color = 13551615

$fcA7 = $ws.Range("A7").FormatConditions.AddUniqueValues()
$fcA7.DupeUnique = 1
$fcA7.Font.Color = 393372
$fcA7.Interior.Color = 13551615

$fcB7 = $ws.Range("B7").FormatConditions.AddUniqueValues()
$fcB7.DupeUnique = 1
$fcB7.Font.Color = 393372
$fcB7.Interior.Color = 13551615

$fcA1314.Priority = 4
$fcB1314.Priority = 3
$fcA7.Priority = 2
$fcB7.Priority = 1

# Selection moved to Q20
$ws.Range("Q20").Select()
